$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hidden")

$rows = @(
    @(2, 'address', 0, 'Adres', $null)
    @(3, 'changePassword', 0, 'Wijzig wachtwoord', 'Changer mot de passe')
    @(4, 'city', 0, 'Stad', $null)
    @(5, 'country', 0, 'Land', $null)
    @(6, 'email', 0, 'E-mail', 'E-mail')
    @(7, 'failed', 0, 'De inloggegevens klopten niet.', $null)
    @(8, 'firstName', 0, 'Voornaam', $null)
    @(9, 'forgotPassword', 0, 'Wachtwoord vergeten?', 'Mot de passe oublié?')
    @(10, 'invalidResetLink', 0, 'Deze link is ongeldig.', 'Ce lien n''est pas valable.')
    @(11, 'lastName', 0, 'Naam', $null)
    @(12, 'loggedOut', 0, 'Je bent nu uitgelogd.', 'Vous êtes déconnecté(e).')
    @(13, 'login', 0, 'Log in', 'Login')
    @(14, 'logout', 0, 'Uitloggen', $null)
    @(15, 'logout.title', 0, 'Logout', 'Logout')
    @(16, 'noAccount', 0, 'Nog geen account?', 'Les données ne sont pas correctes.')
    @(17, 'notActivatedError', 0, 'Uw account is nog niet actief.', 'Votre compte n''est pas encore actif.')
    @(18, 'password', 0, 'Wachtwoord', $null)
    @(19, 'passwordChanged', 0, 'Uw wachtwoord werd gewijzigd.', 'Votre mot de passe a été changé.')
    @(20, 'passwordConfirm', 0, 'Bevestig wachtwoord', 'Confirmer mot de passe')
    @(21, 'postal', 0, 'Postcode', $null)
    @(22, 'register', 0, 'Registreer', $null)
    @(23, 'register.submit', 0, 'Maak profiel aan', 'Créez votre profil')
    @(24, 'register.toLogin', 0, 'Naar login', 'Vers login')
    @(25, 'resetInstructions', 0, 'Je nieuwe wachtwoord moet minstens 8 karakters lang zijn.', 'Votre nouveau mot de passe doit contenir au moins 8 caractères.')
    @(26, 'resetLinkExpired', 0, 'Deze link is niet meer geldig.', 'Ce lien n''est plus valable.')
    @(27, 'resetPassword.button', 0, 'Mail me', 'Envoyez-moi l''e-mail')
    @(28, 'resetPassword.intro', 0, 'Geef je e-mailadres op en we sturen je een link waarmee je je wachtwoord kan wijzigen', 'Donnez-nous votre adresse e-mail et nous vous envoyons un lien par lequel vous pouvez changer votre mot de passe')
    @(29, 'resetPassword.title', 0, 'Wachtwoord opvragen', 'Demander mot de passe')
    @(30, 'resetPassword.toLogin', 0, 'Naar login', 'Vers login')
    @(31, 'resetPasswordButton', 0, 'Mail me', 'Envoyez-moi l''e-mail')
    @(32, 'resetPasswordIntro', 0, 'Geef je e-mailadres op en we sturen je een link waarmee je je wachtwoord kan wijzigen.', 'Donnez votre adresse e-mail et nous vous envoyons un lien par lequel vous pouvez changer votre mot de passe.')
    @(33, 'telephone', 0, 'Telefoon', 'Initialiser mot de passe')
    @(34, 'titleChangePassword', 0, 'Wachtwoord instellen', $null)
    @(35, 'titleLogin', 0, 'Login', 'Login')
    @(36, 'titleRegister', 0, 'Maak profiel aan', 'Créez votre profil')
    @(37, 'titleResetPassword', 0, 'Wachtwoord opvragen', 'Demander mot de passe')
    @(38, 'toLogin', 0, 'Inloggen?', 'Connecter?')
    @(39, 'toRegistrationForm', 0, 'Nog geen profiel?', 'Pas encore de profil?')
    @(40, 'unknownUser', 0, 'Er is niemand geregistreerd met dit e-mailadres.', 'Cette adresse e-mail n''est pas connue.')
    @(41, 'waitingForApproval', 0, 'Uw aanvraag is ontvangen. Er worden er bevestiging gestuurd zodra deze wordt goedgekeurd.', 'Nous avons reçu votre demande. Vous recevrez une confirmation dès qu''elle sera approuvée.')
)

foreach ($row in $rows) {
    $rn = $row[0]
    $ws.Cells.Item($rn, 1).Value = "auth"
    $ws.Cells.Item($rn, 2).Value = $row[1]
    $ws.Cells.Item($rn, 3).Value = $row[2]
    $ws.Cells.Item($rn, 5).Value = $row[3]
    if ($row[4] -ne $null) {
        $ws.Cells.Item($rn, 6).Value = $row[4]
    }
}

$ws.Activate()
$ws.Range("E159").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 2
